$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the activity info: week C4 (01.09) had no activity, but actually hosts "Oppgaveseminar 1",
# while D3 (28.08) was wrongly listed as "Oppgaveseminar 1" but should say "Ingen aktivitet".
$ws.Range("C4").Value = "01.09: **Oppgaveseminar 1** i Aud A. Se \@ref(seminar) for oppgaver."
$ws.Range("D3").Value = "28.08: (Ingen aktivitet)"

# Update the active selection to D4, as recorded in the saved workbook view.
$ws.Range("D4").Select()
